$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing index=1 user) gets new data
$ws.Cells.Item(2, 2).Value = "Lucas"
$ws.Cells.Item(2, 3).Value = 998332344
$ws.Cells.Item(2, 4).Value = "cariolanodcosta@gmail.com"

# Row 3: the original "kaique" row, shifted down, index=2
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "kaique"
$ws.Cells.Item(3, 3).Value = 123
$ws.Cells.Item(3, 4).Value = "kaskjdqo@jnvawf.com"

# Row 4: new user "kaique luan", index=3, password stored as text "123"
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "kaique luan"
$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "123"
$ws.Cells.Item(4, 4).Value = "cariolanodcosta@gmail.com"
